$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Move the selection on "products" from O5 to B1 (view-state tweak from
#    the commit).
# ---------------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("products")
$wsProducts.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Add the two new KPI lookup sheets at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$kpiLevel2 = $wb.Worksheets.Add($null, $lastSheet)
$kpiLevel2.Name = "kpi_level_2"
$kpiResultValues = $wb.Worksheets.Add($null, $kpiLevel2)
$kpiResultValues.Name = "kpi_result_values"

function Write-Row($sheet, $rowNum, $values) {
    for ($i = 0; $i -lt $values.Count; $i++) {
        $v = $values[$i]
        if ($null -ne $v -and $v -ne "") {
            $sheet.Cells.Item($rowNum, $i + 1).Value2 = $v
        }
    }
}

# --- kpi_level_2 --------------------------------------------------------
$kpiLevel2Header = @("pk","kpi_level_1_fk","type","client_name","kpi_family_fk","version", `
    "numerator_type_fk","denominator_type_fk","kpi_score_type_fk","kpi_result_type_fk", `
    "valid_from","valid_until","delete_time","initiated_by","context_type_fk", `
    "kpi_calculation_stage_fk","session_relevance","scene_relevance","planogram_relevance", `
    "live_session_relevance","live_scene_relevance","kpi_target_type_fk","is_percent","kpi_class_path")
Write-Row $kpiLevel2 1 $kpiLevel2Header

Write-Row $kpiLevel2 2 @(1, $null, "TOPSKU_0", "OSA Score", 20, 1, 3, 5, $null, $null, $null, $null, $null, "Custom", 8, 3, 1, 0, 0, 0, 0, $null, 0, $null)
Write-Row $kpiLevel2 3 @(2, $null, "TOPSKU_CATEGORY", "OSA Category", 20, 1, 4, 8, $null, $null, $null, $null, $null, "Custom", $null, 3, 1, 0, 0, 0, 0, $null, 0, $null)
Write-Row $kpiLevel2 4 @(3, $null, "TOPSKU_BUNDLE", "OSA Bundle", 20, 1, 2002, 4, $null, 1, $null, $null, $null, "Custom", 2001, 3, 1, 0, 0, 0, 0, $null, 0, $null)
Write-Row $kpiLevel2 5 @(4, $null, "TOPSKU_SKU", "OSA SKU", 20, 1, 1, 2002, $null, 1, $null, $null, $null, "Custom", $null, 3, 1, 0, 0, 0, 0, $null, 0, $null)

# header formatting: copy the existing bold/border/center-top header style
# used throughout the workbook (e.g. products!A1) onto row 1.
$wsProducts.Range("A1").Copy() | Out-Null
$kpiLevel2.Range("A1:X1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- kpi_result_values ---------------------------------------------------
Write-Row $kpiResultValues 1 @("pk","value","kpi_result_type_fk")
Write-Row $kpiResultValues 2 @(1, "OOS", 1)
Write-Row $kpiResultValues 3 @(2, "DISTRIBUTED", 1)
Write-Row $kpiResultValues 4 @(3, "EXTRA", 1)

$wsProducts.Range("A1").Copy() | Out-Null
$kpiResultValues.Range("A1:C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------------
# 3. New autofilter-derived defined name on "scif" (LibreOffice-style
#    incrementing _xlnm._FilterDatabase_* entry added on save).
# ---------------------------------------------------------------------------
$wsScif = $wb.Worksheets.Item("scif")
$wsScif.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0_0", "=scif!`$A`$1:`$Z`$36") | Out-Null

# ---------------------------------------------------------------------------
# 4. Make "kpi_level_2" the active sheet/tab (activeTab goes from 2 to 5)
#    with the cursor sitting on A7, matching the saved view state.
# ---------------------------------------------------------------------------
$kpiLevel2.Activate()
$kpiLevel2.Range("A7").Select() | Out-Null

$kpiResultValues.Range("C2").Select() | Out-Null
$kpiLevel2.Activate()
$kpiLevel2.Range("A7").Select() | Out-Null

Write-Output "edit complete"
